# Pooh Points: normal 20260214
# - Games that were in-progress ("<clock> - 2nd Half") have finished; mark their
#   status as "Final".
# - Two players' minutes (P53, P117) were corrected down by 1 once the final
#   box score posted.
# - The "status" column (G) no longer needs to be as wide now that every row
#   reads "Final" instead of a running clock string, so narrow it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Rows whose in-progress game clock has become "Final".
$finishedRows = @(
    2, 5, 6, 7, 12, 13, 15, 17, 20, 26, 28, 32, 35, 40, 41, 42, 43, 45, 48,
    53, 54, 55, 61, 62, 67, 68, 71, 73, 75, 80, 85, 88, 94, 97, 105, 112,
    113, 114, 117, 118, 119, 120, 121, 123, 126, 129, 130, 136, 137, 138, 139
)

foreach ($row in $finishedRows) {
    $ws.Range("G$row").Value = "Final"
}

# Final box score correction: minutes played.
$ws.Range("P53").Value = 25
$ws.Range("P117").Value = 9

# Shrink the now-shorter "status" column.
$ws.Columns.Item(7).ColumnWidth = 7.17
